# Apply updated cryptocurrency price/volume data to Sheet1.
# Numeric-looking "Price" values are written with a leading apostrophe so that
# Excel stores them as text (matching the workbook's original inline-string/text cells)
# instead of silently converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '43.773.94'
$ws.Cells.Item(2, 5).Value = '  -0.06%  '
# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '2.282.20'
$ws.Cells.Item(3, 5).Value = '  -0.46%  '
# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = '  -0.45%  '
# Row 5: Solana
$ws.Cells.Item(5, 4).Value = '''123.85'
$ws.Cells.Item(5, 5).Value = '  +9.14%  '
# Row 6: BNB
$ws.Cells.Item(6, 4).Value = '''266.17'
$ws.Cells.Item(6, 5).Value = '  -0.94%  '
# Row 7: XRP
$ws.Cells.Item(7, 4).Value = '''0.637'
$ws.Cells.Item(7, 5).Value = '  +1.73%  '
# Row 8: USDC
$ws.Cells.Item(8, 5).Value = '  +0.17%  '
# Row 9: Cardano
$ws.Cells.Item(9, 4).Value = '''0.625'
$ws.Cells.Item(9, 5).Value = '  +1.35%  '
# Row 10: Avalanche
$ws.Cells.Item(10, 4).Value = '''48.57'
$ws.Cells.Item(10, 5).Value = '  +0.43%  '
# Row 11: Dogecoin
$ws.Cells.Item(11, 5).Value = '  +0.49%  '
# Row 12: Polkadot
$ws.Cells.Item(12, 4).Value = '''9.15'
$ws.Cells.Item(12, 5).Value = '  +1.70%  '
# Row 13: TRON
$ws.Cells.Item(13, 5).Value = '  -1.01%  '
# Row 14: Chainlink
$ws.Cells.Item(14, 4).Value = '''15.48'
$ws.Cells.Item(14, 5).Value = '  -1.95%  '
# Row 15: Polygon
$ws.Cells.Item(15, 4).Value = '''0.900'
$ws.Cells.Item(15, 5).Value = '  +4.59%  '
# Row 16: WrappedliquidstakedEther2.0
$ws.Cells.Item(16, 4).Value = '2.624.18'
$ws.Cells.Item(16, 5).Value = '  -0.54%  '
# Row 17: WrappedEther
$ws.Cells.Item(17, 4).Value = '2.280.45'
$ws.Cells.Item(17, 5).Value = '  -0.54%  '
# Row 18: WrappedBTC
$ws.Cells.Item(18, 4).Value = '43.718.16'
$ws.Cells.Item(18, 5).Value = '  +0.06%  '
# Row 19: ShibaInu
$ws.Cells.Item(19, 5).Value = '  +0.29%  '
# Row 20: Uniswap
$ws.Cells.Item(20, 5).Value = '  -0.42%  '
# Row 21: Litecoin
$ws.Cells.Item(21, 4).Value = '''72.42'
$ws.Cells.Item(21, 5).Value = '  +0.18%  '
# Row 22: ImmutableX
$ws.Cells.Item(22, 5).Value = '  +0.83%  '
# Row 23: BitcoinCash
$ws.Cells.Item(23, 4).Value = '''238.07'
$ws.Cells.Item(23, 5).Value = '  +2.19%  '
# Row 24: InternetComputer(DFINITY)
$ws.Cells.Item(24, 4).Value = '''9.54'
$ws.Cells.Item(24, 5).Value = '  -3.08%  '
# Row 25: PancakeSwap
$ws.Cells.Item(25, 4).Value = '''2.89'
$ws.Cells.Item(25, 5).Value = '  -1.65%  '
# Row 26: Dai
$ws.Cells.Item(26, 5).Value = '  +1.87%  '
# Row 27: Cosmos
$ws.Cells.Item(27, 4).Value = '''11.87'
$ws.Cells.Item(27, 5).Value = '  +0.99%  '
# Row 28: InjectiveProtocol
$ws.Cells.Item(28, 4).Value = '''42.87'
$ws.Cells.Item(28, 5).Value = '  -1.63%  '
# Row 29: WEMIXToken
$ws.Cells.Item(29, 5).Value = '  -0.48%  '
# Row 30: Toncoin
$ws.Cells.Item(30, 5).Value = '  +0.13%  '
# Row 31: Monero
$ws.Cells.Item(31, 4).Value = '''172.34'
$ws.Cells.Item(31, 5).Value = '  -1.51%  '
# Row 32: EthereumClassic
$ws.Cells.Item(32, 4).Value = '''21.77'
$ws.Cells.Item(32, 5).Value = '  +0.64%  '
# Row 33: Hedera
$ws.Cells.Item(33, 4).Value = '''0.0913'
$ws.Cells.Item(33, 5).Value = '  -2.03%  '
# Row 34: Filecoin
$ws.Cells.Item(34, 4).Value = '''5.75'
$ws.Cells.Item(34, 5).Value = '  +1.32%  '
# Row 35: Stellar
$ws.Cells.Item(35, 5).Value = '  +1.96%  '
# Row 36: NEARProtocol
$ws.Cells.Item(36, 4).Value = '''4.21'
$ws.Cells.Item(36, 5).Value = '  +10.73%  '
# Row 37: VeChain
$ws.Cells.Item(37, 5).Value = '  +5.38%  '
# Row 38: RenderToken
$ws.Cells.Item(38, 5).Value = '  -2.63%  '
# Row 39: Kaspa
$ws.Cells.Item(39, 5).Value = '  +1.31%  '
# Row 40: LidoDAOToken
$ws.Cells.Item(40, 4).Value = '''2.54'
$ws.Cells.Item(40, 5).Value = '  +5.73%  '
# Row 41: MultiversX
$ws.Cells.Item(41, 4).Value = '''75.33'
$ws.Cells.Item(41, 5).Value = '  -0.11%  '
# Row 42: Celestia
$ws.Cells.Item(42, 5).Value = '  -2.46%  '
# Row 43: Algorand
$ws.Cells.Item(43, 4).Value = '''0.240'
$ws.Cells.Item(43, 5).Value = '  -0.77%  '
# Row 44: FirstDigitalUSD
$ws.Cells.Item(44, 4).Value = '''0.999'
$ws.Cells.Item(44, 5).Value = '  -0.43%  '
# Row 45: ARBITRUM
$ws.Cells.Item(45, 5).Value = '  -2.89%  '
# Row 46: THORChain
$ws.Cells.Item(46, 4).Value = '''5.69'
$ws.Cells.Item(46, 5).Value = '  -10.40%  '
# Row 47: ordi
$ws.Cells.Item(47, 4).Value = '''74.29'
$ws.Cells.Item(47, 5).Value = '  +37.27%  '
# Row 48: TrustWalletToken -> FraxShare
$ws.Cells.Item(48, 2).Value = 'FraxShare'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(48, 4).Value = '''8.62'
$ws.Cells.Item(48, 5).Value = '  -1.99%  '
# Row 49: FraxShare -> TrustWalletToken
$ws.Cells.Item(49, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(49, 4).Value = '''1.27'
$ws.Cells.Item(49, 5).Value = '  +0.31%  '
# Row 50: Cronos
$ws.Cells.Item(50, 4).Value = '''0.101'
$ws.Cells.Item(50, 5).Value = '  +1.09%  '
# Row 51: Aave
$ws.Cells.Item(51, 4).Value = '''102.35'
$ws.Cells.Item(51, 5).Value = '  -0.58%  '
